$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeroDocumento column (B) and usuario column (D) for rows 2-5
foreach ($r in 2..5) {
    $ws.Cells.Item($r, 2).Value = "22483228"
    $ws.Cells.Item($r, 4).Value = "userrobot2"
}

# Row 3: account type/number/amount mapping changes
$ws.Range("P3").Value = "Corriente"
$ws.Range("Q3").Value = "406-132280-02"
$ws.Range("R3").Value = "150000"

# Update the view / selection to match the saved workbook state
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("N9").Select()
